$d = $word.ActiveDocument

# 1) "10:45am  Intro to Exploratory Data Analysis (1h 15min) - Marcus Beck"
#    becomes "10:45am  Data Processing (1h 15min) - Marcus Beck"
$d.Content.Find.Execute("Intro to Exploratory Data Analysis", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Data Processing", 2)

# 2) "01:00pm  Exploratory Data Analysis - Continued (1 hr) - Marcus Beck"
#    becomes "01:00pm  Exploratory Data Analysis (1 hr) - Marcus Beck"
$d.Content.Find.Execute("Exploratory Data Analysis - Continued", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Exploratory Data Analysis", 2)
